$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# windMun.* columns (J1:T1)
$ws.Range("J1").Value = "windMun.42"
$ws.Range("K1").Value = "windMun.34"
$ws.Range("L1").Value = "windMun.15"
$ws.Range("M1").Value = "windMun.18"
$ws.Range("N1").Value = "windMun.03"
$ws.Range("O1").Value = "windMun.11"
$ws.Range("P1").Value = "windMun.38"
$ws.Range("Q1").Value = "windMun.54"
$ws.Range("R1").Value = "windMun.50"
$ws.Range("S1").Value = "windMun.46"
$ws.Range("T1").Value = "windMun.30"

# sunPower.* columns (AK1:AU1)
$ws.Range("AK1").Value = "sunPower.42"
$ws.Range("AL1").Value = "sunPower.34"
$ws.Range("AM1").Value = "sunPower.15"
$ws.Range("AN1").Value = "sunPower.18"
$ws.Range("AO1").Value = "sunPower.03"
$ws.Range("AP1").Value = "sunPower.11"
$ws.Range("AQ1").Value = "sunPower.38"
$ws.Range("AR1").Value = "sunPower.54"
$ws.Range("AS1").Value = "sunPower.50"
$ws.Range("AT1").Value = "sunPower.46"
$ws.Range("AU1").Value = "sunPower.30"
